$d = $word.ActiveDocument

# Position an insertion point at the very end of the document's main story.
$contentEnd = $d.Content.End
$r = $d.Range($contentEnd, $contentEnd)

# Build the OOXML fragment for the two new paragraphs:
#   1) a blank paragraph with no formatting
#   2) a paragraph with the "My thoughts for you..." note, including the
#      w:proofErr gramStart/gramEnd markers around "averaged" and "raters"
#      exactly as they appear in the target revision.
$body = ""
$body += '<w:p/>'
$body += '<w:p>'
$body += '<w:r><w:t xml:space="preserve">My thoughts for you to cut and paste for later: teach a little about the methods including explaining an interpretation of an odds ratio for the different models. Make a Table 1 </w:t></w:r>'
$body += '<w:proofErr w:type="gramStart"/>'
$body += '<w:r><w:t>averaged</w:t></w:r>'
$body += '<w:proofErr w:type="gramEnd"/>'
$body += '<w:r><w:t xml:space="preserve"> across </w:t></w:r>'
$body += '<w:proofErr w:type="gramStart"/>'
$body += '<w:r><w:t>raters</w:t></w:r>'
$body += '<w:proofErr w:type="gramEnd"/>'
$body += '<w:r><w:t xml:space="preserve"> or make some graphs of the raw data on percentages in each severity.  Include other measures of agreement from raters.</w:t></w:r>'
$body += '</w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $body + '</w:body></w:document></pkg:xmlData>' + `
  '</pkg:part>' + `
'</pkg:package>'

$r.InsertXML($xml)
